# Updates the cryptocurrency price/volume columns (and, for rows 41-44,
# the coin name/link) to match the refreshed GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextCell($ws, $addr, $val) {
    # Leading apostrophe forces text entry (mirrors manual typing),
    # preventing Excel from re-parsing numeric-looking strings (e.g.
    # "0.380", "37.348.10") and mangling their literal formatting.
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

Set-TextCell $ws 'D2' '37.348.10'
Set-TextCell $ws 'E2' '  +2.02%  '

Set-TextCell $ws 'D3' '2.090.48'
Set-TextCell $ws 'E3' '  +0.35%  '

Set-TextCell $ws 'E4' '  +0.12%  '

Set-TextCell $ws 'D5' '251.88'
Set-TextCell $ws 'E5' '  +1.72%  '

Set-TextCell $ws 'D6' '0.666'
Set-TextCell $ws 'E6' '  +0.38%  '

Set-TextCell $ws 'E7' '  -0.06%  '

Set-TextCell $ws 'D8' '54.10'
Set-TextCell $ws 'E8' '  +19.78%  '

Set-TextCell $ws 'D9' '62.45'
Set-TextCell $ws 'E9' '  +2.62%  '

Set-TextCell $ws 'D10' '0.380'
Set-TextCell $ws 'E10' '  +4.41%  '

Set-TextCell $ws 'D11' '0.0753'
Set-TextCell $ws 'E11' '  +4.09%  '

Set-TextCell $ws 'E12' '  +7.70%  '

Set-TextCell $ws 'D13' '15.45'
Set-TextCell $ws 'E13' '  +6.56%  '

Set-TextCell $ws 'D14' '2.396.35'
Set-TextCell $ws 'E14' '  +0.57%  '

Set-TextCell $ws 'D15' '0.852'
Set-TextCell $ws 'E15' '  +2.74%  '

Set-TextCell $ws 'D16' '5.28'
Set-TextCell $ws 'E16' '  +7.36%  '

Set-TextCell $ws 'D17' '2.090.91'
Set-TextCell $ws 'E17' '  +0.41%  '

Set-TextCell $ws 'D18' '37.311.16'
Set-TextCell $ws 'E18' '  +1.69%  '

Set-TextCell $ws 'D19' '73.17'
Set-TextCell $ws 'E19' '  +2.20%  '

Set-TextCell $ws 'D20' '14.45'
Set-TextCell $ws 'E20' '  +13.53%  '

Set-TextCell $ws 'E21' '  +4.72%  '

Set-TextCell $ws 'D22' '241.40'
Set-TextCell $ws 'E22' '  +1.09%  '

Set-TextCell $ws 'D23' '5.28'
Set-TextCell $ws 'E23' '  +7.21%  '

Set-TextCell $ws 'E24' '  -0.18%  '

Set-TextCell $ws 'D25' '2.48'
Set-TextCell $ws 'E25' '  +0.58%  '

Set-TextCell $ws 'D26' '172.07'
Set-TextCell $ws 'E26' '  +1.60%  '

Set-TextCell $ws 'E27' '  +4.75%  '

Set-TextCell $ws 'D28' '20.94'
Set-TextCell $ws 'E28' '  +2.86%  '

Set-TextCell $ws 'E29' '  +3.54%  '

Set-TextCell $ws 'E30' '  +2.49%  '

Set-TextCell $ws 'D31' '23.65'
Set-TextCell $ws 'E31' '  +8.50%  '

Set-TextCell $ws 'D32' '1.11'
Set-TextCell $ws 'E32' '  +22.77%  '

Set-TextCell $ws 'D33' '4.55'
Set-TextCell $ws 'E33' '  +4.45%  '

Set-TextCell $ws 'D34' '0.0627'
Set-TextCell $ws 'E34' '  +7.42%  '

Set-TextCell $ws 'D35' '0.0902'
Set-TextCell $ws 'E35' '  -0.62%  '

Set-TextCell $ws 'D36' '4.29'
Set-TextCell $ws 'E36' '  +7.20%  '

Set-TextCell $ws 'E37' '  -0.08%  '

Set-TextCell $ws 'D38' '2.27'
Set-TextCell $ws 'E38' '  -1.17%  '

Set-TextCell $ws 'D39' '1.83'
Set-TextCell $ws 'E39' '  -3.78%  '

Set-TextCell $ws 'E40' '  +1.95%  '

Set-TextCell $ws 'B41' 'VeChain'
Set-TextCell $ws 'C41' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell $ws 'D41' '0.0229'
Set-TextCell $ws 'E41' '  +6.25%  '

Set-TextCell $ws 'B42' 'InjectiveProtocol'
Set-TextCell $ws 'C42' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextCell $ws 'D42' '17.91'
Set-TextCell $ws 'E42' '  +12.41%  '

Set-TextCell $ws 'B43' 'ARBITRUM'
Set-TextCell $ws 'C43' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell $ws 'D43' '1.18'
Set-TextCell $ws 'E43' '  +3.65%  '

Set-TextCell $ws 'B44' 'FTXToken'
Set-TextCell $ws 'C44' 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextCell $ws 'D44' '4.56'
Set-TextCell $ws 'E44' '  +129.48%  '

Set-TextCell $ws 'D45' '99.65'
Set-TextCell $ws 'E45' '  +2.05%  '

Set-TextCell $ws 'D46' '0.0970'
Set-TextCell $ws 'E46' '  +18.09%  '

Set-TextCell $ws 'D47' '2.81'
Set-TextCell $ws 'E47' '  +0.37%  '

Set-TextCell $ws 'D48' '1.334.20'
Set-TextCell $ws 'E48' '  +0.45%  '

Set-TextCell $ws 'E49' '  +3.59%  '

Set-TextCell $ws 'E50' '  +6.68%  '

Set-TextCell $ws 'D51' '6.97'
Set-TextCell $ws 'E51' '  +12.69%  '

